$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1377
$ws.Range("I4").Value = 1543.1428
$ws.Range("K4").Value = 1543.1428
$ws.Range("M4").Value = -1429.1428

$ws.Range("H5").Value = 33.75
$ws.Range("I5").Value = 33.75
$ws.Range("K5").Value = 33.75
$ws.Range("M5").Value = 81.25

$ws.Range("H40").Value = 3939.8
$ws.Range("I40").Value = 1700
$ws.Range("J40").Value = 4499.75
$ws.Range("K40").Value = 1700
$ws.Range("L40").Value = 4499.75
$ws.Range("M40").Value = -1525
$ws.Range("N40").Value = -4849.75

$ws.Range("H53").Value = 4615.6
$ws.Range("I53").Value = 8444.6
$ws.Range("K53").Value = 8444.6
$ws.Range("M53").Value = -7807.6

$ws.Range("H112").Value = 1712.6923
$ws.Range("J112").Value = 1790.6666
$ws.Range("L112").Value = 5371.9998
$ws.Range("N112").Value = -7587.9998

$ws.Range("H136").Value = 68666.664
$ws.Range("J136").Value = 68666.664
$ws.Range("L136").Value = 68666.664
$ws.Range("N136").Value = -78866.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 500
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -732

$ws.Range("H32").Value = 6848.7026
$ws.Range("I32").Value = 6729.648
$ws.Range("K32").Value = 6729.648
$ws.Range("M32").Value = -6442.648

$ws.Range("H44").Value = 88000
$ws.Range("J44").Value = 88000
$ws.Range("L44").Value = 88000
$ws.Range("N44").Value = -88976

$ws.Range("H61").Value = 10012.479
$ws.Range("I61").Value = 11405.444
$ws.Range("K61").Value = 11405.444
$ws.Range("M61").Value = -11193.444

$ws.Range("H74").Value = 5233.5864
$ws.Range("J74").Value = 2327.7222
$ws.Range("L74").Value = 2327.7222
$ws.Range("N74").Value = -4075.7222

$ws.Range("H77").Value = 5233.5864
$ws.Range("J77").Value = 2327.7222
$ws.Range("L77").Value = 11638.611
$ws.Range("N77").Value = -20374.611

$ws.Range("H113").Value = 116999.555
$ws.Range("J113").Value = 122874.5
$ws.Range("L113").Value = 122874.5
$ws.Range("N113").Value = -131552.5

$ws.Range("H122").Value = 1039740.94
$ws.Range("I122").Value = 4978.125
$ws.Range("J122").Value = 6006602.5
$ws.Range("K122").Value = 14934.375
$ws.Range("L122").Value = 18019807.5
$ws.Range("M122").Value = -12484.375
$ws.Range("N122").Value = -18024707.5

$ws.Range("H132").Value = 3244.577
$ws.Range("I132").Value = 2359.6875
$ws.Range("J132").Value = 4660.4
$ws.Range("K132").Value = 7079.0625
$ws.Range("L132").Value = 13981.2
$ws.Range("M132").Value = -4549.0625
$ws.Range("N132").Value = -19041.2

$ws.Range("H136").Value = 10012.479
$ws.Range("I136").Value = 11405.444
$ws.Range("K136").Value = 34216.33199999999
$ws.Range("M136").Value = -31666.33199999999

$ws.Range("H141").Value = 43627.875
$ws.Range("J141").Value = 43627.875
$ws.Range("L141").Value = 43627.875
$ws.Range("N141").Value = -53987.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 70000
$ws.Range("J92").Value = 70000
$ws.Range("L92").Value = 70000
$ws.Range("N92").Value = -74992

$ws.Range("H134").Value = 7819.16
$ws.Range("I134").Value = 8972.611000000001
$ws.Range("K134").Value = 26917.833
$ws.Range("M134").Value = -24382.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6710.857
$ws.Range("I31").Value = 6746
$ws.Range("J31").Value = 6500
$ws.Range("K31").Value = 6746
$ws.Range("L31").Value = 6500
$ws.Range("M31").Value = -6451
$ws.Range("N31").Value = -7090

$ws.Range("H34").Value = 6710.857
$ws.Range("I34").Value = 6746
$ws.Range("J34").Value = 6500
$ws.Range("K34").Value = 6746
$ws.Range("L34").Value = 6500
$ws.Range("M34").Value = -6544
$ws.Range("N34").Value = -6904

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H58").Value = 3530.24
$ws.Range("I58").Value = 3920.75
$ws.Range("K58").Value = 3920.75
$ws.Range("M58").Value = -3717.75

$ws.Range("H136").Value = 3530.24
$ws.Range("I136").Value = 3920.75
$ws.Range("K136").Value = 11762.25
$ws.Range("M136").Value = -9212.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 69405.62
$ws.Range("J80").Value = 74422.836
$ws.Range("L80").Value = 223268.508
$ws.Range("N80").Value = -225140.508

$ws.Range("H83").Value = 69405.62
$ws.Range("J83").Value = 74422.836
$ws.Range("L83").Value = 669805.524
$ws.Range("N83").Value = -679165.524

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 287.7143
$ws.Range("I2").Value = 288.8
$ws.Range("J2").Value = 285
$ws.Range("K2").Value = 288.8
$ws.Range("L2").Value = 285
$ws.Range("M2").Value = -175.8
$ws.Range("N2").Value = -511

$ws.Range("H53").Value = 20000
$ws.Range("I53").Value = 20000
$ws.Range("K53").Value = 20000
$ws.Range("M53").Value = -19369

$ws.Range("H119").Value = 28000
$ws.Range("J119").Value = 28000
$ws.Range("L119").Value = 28000
$ws.Range("N119").Value = -37676

$ws.Range("H123").Value = 22065.2
$ws.Range("J123").Value = 22065.2
$ws.Range("L123").Value = 22065.2
$ws.Range("N123").Value = -26965.2

$ws.Range("H126").Value = 6991.1724
$ws.Range("I126").Value = 9394
$ws.Range("J126").Value = 3587.1667
$ws.Range("K126").Value = 28182
$ws.Range("L126").Value = 10761.5001
$ws.Range("M126").Value = -25712
$ws.Range("N126").Value = -15701.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 499
$ws.Range("J9").Value = 466.66666
$ws.Range("L9").Value = 466.66666
$ws.Range("N9").Value = -914.66666

$ws.Range("H55").Value = 1234.6471
$ws.Range("I55").Value = 337.9091
$ws.Range("J55").Value = 2878.6667
$ws.Range("K55").Value = 337.9091
$ws.Range("L55").Value = 2878.6667
$ws.Range("M55").Value = -164.9091
$ws.Range("N55").Value = -3224.6667

$ws.Range("H61").Value = 4878
$ws.Range("J61").Value = 12899.6
$ws.Range("L61").Value = 12899.6
$ws.Range("N61").Value = -13303.6

$ws.Range("H74").Value = 41998.5
$ws.Range("I74").Value = 41998.5
$ws.Range("K74").Value = 41998.5
$ws.Range("M74").Value = -41000.5

$ws.Range("H77").Value = 41998.5
$ws.Range("I77").Value = 41998.5
$ws.Range("K77").Value = 125995.5
$ws.Range("M77").Value = -121003.5

$ws.Range("H113").Value = 4878
$ws.Range("J113").Value = 12899.6
$ws.Range("L113").Value = 12899.6
$ws.Range("N113").Value = -17239.6

$ws.Range("H122").Value = 4821.5
$ws.Range("I122").Value = 5024.684
$ws.Range("K122").Value = 15074.052
$ws.Range("M122").Value = -12624.052

$ws.Range("H136").Value = 5362.316
$ws.Range("I136").Value = 4173.375
$ws.Range("K136").Value = 12520.125
$ws.Range("M136").Value = -9970.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 59954.5
$ws.Range("J133").Value = 59954.5
$ws.Range("L133").Value = 59954.5
$ws.Range("N133").Value = -70074.5

$ws.Range("H136").Value = 347460.66
$ws.Range("I136").Value = 431330
$ws.Range("K136").Value = 1293990
$ws.Range("M136").Value = -1291440

$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
